# ornek_gider.xlsx — "Add files via upload" edit
#
# - A2:A4  "ETKİ AKADEMİ"   -> "Etki Akademi"
# - C2     "Kira"            -> "Araç Sigortası"
# - C3     "Araç Sigortası"  -> "İşyeri Kirası"
# - C4     "Elektrik"        -> "Yıllık Yazılım"
# - D3     24000             -> 22750
# - D4     1800              -> 24000
# - F2     2025-06-01        -> 2025-06-30
# - F3     2026-05-31        -> 2025-06-30
# - F4     2025-06-30        -> 2026-05-31
#
# Columns B, E and F hold dates typed in as plain text (General format,
# shared-string cells), not real Excel date serials. Assigning a
# "yyyy-mm-dd"-shaped literal straight to .Value lets Excel's smart
# entry reinterpret it as a date, which would both change the stored
# type and stamp a date NumberFormat on the cell. Prefixing the literal
# with an apostrophe forces text entry instead; that alone marks the
# cell "quote prefixed" (a distinct style), so the style is reset back
# to Normal right after the value is in, leaving the cell exactly as
# plain/unstyled as it started.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Etki Akademi"
$ws.Range("A3").Value = "Etki Akademi"
$ws.Range("A4").Value = "Etki Akademi"

$ws.Range("C2").Value = "Araç Sigortası"
$ws.Range("C3").Value = "İşyeri Kirası"
$ws.Range("C4").Value = "Yıllık Yazılım"

$ws.Range("D3").Value = 22750
$ws.Range("D4").Value = 24000

$ws.Range("F2").Value = "'2025-06-30"
$ws.Range("F2").Style = "Normal"

$ws.Range("F3").Value = "'2025-06-30"
$ws.Range("F3").Style = "Normal"

$ws.Range("F4").Value = "'2026-05-31"
$ws.Range("F4").Style = "Normal"
